$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216
$ws.Range("H40").Value = 22620330
$ws.Range("J40").Value = 100000000
$ws.Range("L40").Value = 100000000
$ws.Range("N40").Value = -100000350
$ws.Range("H58").Value = 11153.625
$ws.Range("I58").Value = 2076.6667
$ws.Range("K58").Value = 6230.000100000001
$ws.Range("M58").Value = -6080.000100000001
$ws.Range("H87").Value = 230000
$ws.Range("J87").Value = 230000
$ws.Range("L87").Value = 230000
$ws.Range("N87").Value = -232496
$ws.Range("H90").Value = 230000
$ws.Range("J90").Value = 230000
$ws.Range("L90").Value = 690000
$ws.Range("N90").Value = -702480
$ws.Range("H92").Value = 1998.25
$ws.Range("I92").Value = 1998.25
$ws.Range("K92").Value = 1998.25
$ws.Range("M92").Value = -750.25
$ws.Range("H98").Value = 6251335
$ws.Range("I98").Value = 7813204.5
$ws.Range("J98").Value = 3858.5
$ws.Range("K98").Value = 7813204.5
$ws.Range("L98").Value = 3858.5
$ws.Range("M98").Value = -7811706.5
$ws.Range("N98").Value = -6854.5
$ws.Range("H100").Value = 3562.24
$ws.Range("I100").Value = 2768.1333
$ws.Range("K100").Value = 2768.1333
$ws.Range("M100").Value = -2227.1333
$ws.Range("H107").Value = 935.93335
$ws.Range("I107").Value = 1048.0769
$ws.Range("J107").Value = 207
$ws.Range("K107").Value = 1048.0769
$ws.Range("L107").Value = 207
$ws.Range("M107").Value = 871.9231
$ws.Range("N107").Value = -4047
$ws.Range("H112").Value = 2366.3333
$ws.Range("J112").Value = 2366.3333
$ws.Range("L112").Value = 7098.999899999999
$ws.Range("N112").Value = -9314.999899999999
$ws.Range("H113").Value = 12598.066
$ws.Range("I113").Value = 11809.25
$ws.Range("J113").Value = 13499.571
$ws.Range("K113").Value = 11809.25
$ws.Range("L113").Value = 13499.571
$ws.Range("M113").Value = -8555.25
$ws.Range("N113").Value = -20007.571
$ws.Range("H121").Value = 2484
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H122").Value = 6251335
$ws.Range("I122").Value = 7813204.5
$ws.Range("J122").Value = 3858.5
$ws.Range("K122").Value = 23439613.5
$ws.Range("L122").Value = 11575.5
$ws.Range("M122").Value = -23437163.5
$ws.Range("N122").Value = -16475.5
$ws.Range("H129").Value = 1319.85
$ws.Range("I129").Value = 816.7143
$ws.Range("J129").Value = 2493.8333
$ws.Range("K129").Value = 2450.1429
$ws.Range("L129").Value = 7481.499899999999
$ws.Range("M129").Value = 2549.8571
$ws.Range("N129").Value = -17481.4999
$ws.Range("H137").Value = 38475650
$ws.Range("I137").Value = 62520890
$ws.Range("J137").Value = 3269.2
$ws.Range("K137").Value = 187562670
$ws.Range("L137").Value = 9807.599999999999
$ws.Range("M137").Value = -187560120
$ws.Range("N137").Value = -14907.6
$ws.Range("H138").Value = 6066.445
$ws.Range("J138").Value = 6066.445
$ws.Range("L138").Value = 18199.335
$ws.Range("N138").Value = -28479.335

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1139.0952
$ws.Range("I2").Value = 1168.12
$ws.Range("J2").Value = 1096.4117
$ws.Range("K2").Value = 1168.12
$ws.Range("L2").Value = 1096.4117
$ws.Range("M2").Value = -1055.12
$ws.Range("N2").Value = -1322.4117
$ws.Range("H32").Value = 19643.115
$ws.Range("I32").Value = 9767.571
$ws.Range("J32").Value = 31164.584
$ws.Range("K32").Value = 9767.571
$ws.Range("L32").Value = 31164.584
$ws.Range("M32").Value = -9480.571
$ws.Range("N32").Value = -31738.584
$ws.Range("H61").Value = 5734.55
$ws.Range("J61").Value = 8577.571
$ws.Range("L61").Value = 8577.571
$ws.Range("N61").Value = -9001.571
$ws.Range("H74").Value = 12616.308
$ws.Range("I74").Value = 8195.6
$ws.Range("J74").Value = 27352
$ws.Range("K74").Value = 8195.6
$ws.Range("L74").Value = 27352
$ws.Range("M74").Value = -7321.6
$ws.Range("N74").Value = -29100
$ws.Range("H77").Value = 12616.308
$ws.Range("I77").Value = 8195.6
$ws.Range("J77").Value = 27352
$ws.Range("K77").Value = 40978
$ws.Range("L77").Value = 136760
$ws.Range("M77").Value = -36610
$ws.Range("N77").Value = -145496
$ws.Range("H110").Value = 4007.2
$ws.Range("I110").Value = 2499.5
$ws.Range("J110").Value = 5012.3335
$ws.Range("K110").Value = 2499.5
$ws.Range("L110").Value = 5012.3335
$ws.Range("M110").Value = -454.5
$ws.Range("N110").Value = -9102.333500000001
$ws.Range("H116").Value = 1139.0952
$ws.Range("I116").Value = 1168.12
$ws.Range("J116").Value = 1096.4117
$ws.Range("K116").Value = 1168.12
$ws.Range("L116").Value = 1096.4117
$ws.Range("M116").Value = 1125.88
$ws.Range("N116").Value = -5684.411700000001
$ws.Range("H122").Value = 5446.375
$ws.Range("I122").Value = 5443.39
$ws.Range("K122").Value = 16330.17
$ws.Range("M122").Value = -13880.17
$ws.Range("H132").Value = 8337.129000000001
$ws.Range("I132").Value = 6418.051
$ws.Range("J132").Value = 19166.215
$ws.Range("K132").Value = 19254.153
$ws.Range("L132").Value = 57498.645
$ws.Range("M132").Value = -16724.153
$ws.Range("N132").Value = -62558.645
$ws.Range("H136").Value = 5734.55
$ws.Range("J136").Value = 8577.571
$ws.Range("L136").Value = 25732.713
$ws.Range("N136").Value = -30832.713

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1139.0952
$ws.Range("I3").Value = 1168.12
$ws.Range("J3").Value = 1096.4117
$ws.Range("K3").Value = 1168.12
$ws.Range("L3").Value = 1096.4117
$ws.Range("M3").Value = -1054.12
$ws.Range("N3").Value = -1324.4117
$ws.Range("H86").Value = 4909.5
$ws.Range("I86").Value = 4565.1665
$ws.Range("J86").Value = 5942.5
$ws.Range("K86").Value = 4565.1665
$ws.Range("L86").Value = 5942.5
$ws.Range("M86").Value = -3442.1665
$ws.Range("N86").Value = -8188.5
$ws.Range("H89").Value = 4909.5
$ws.Range("I89").Value = 4565.1665
$ws.Range("J89").Value = 5942.5
$ws.Range("K89").Value = 22825.8325
$ws.Range("L89").Value = 29712.5
$ws.Range("M89").Value = -17209.8325
$ws.Range("N89").Value = -40944.5
$ws.Range("H94").Value = 2525.5588
$ws.Range("I94").Value = 1541.5238
$ws.Range("J94").Value = 4115.154
$ws.Range("K94").Value = 1541.5238
$ws.Range("L94").Value = 4115.154
$ws.Range("M94").Value = -1090.5238
$ws.Range("N94").Value = -5017.154
$ws.Range("H105").Value = 1701.7142
$ws.Range("I105").Value = 1621.75
$ws.Range("J105").Value = 1957.6
$ws.Range("K105").Value = 1621.75
$ws.Range("L105").Value = 1957.6
$ws.Range("M105").Value = 125.25
$ws.Range("N105").Value = -5451.6

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1878.1
$ws.Range("I16").Value = 1624.3334
$ws.Range("K16").Value = 1624.3334
$ws.Range("M16").Value = -1337.3334
$ws.Range("H20").Value = 89780
$ws.Range("J20").Value = 89780
$ws.Range("L20").Value = 89780
$ws.Range("N20").Value = -90252
$ws.Range("H22").Value = 5623.25
$ws.Range("I22").Value = 3995
$ws.Range("J22").Value = 6600.2
$ws.Range("K22").Value = 3995
$ws.Range("L22").Value = 6600.2
$ws.Range("M22").Value = -3645
$ws.Range("N22").Value = -7300.2
$ws.Range("H30").Value = 89780
$ws.Range("J30").Value = 89780
$ws.Range("L30").Value = 89780
$ws.Range("N30").Value = -89962
$ws.Range("H31").Value = 433127.62
$ws.Range("I31").Value = 9232.75
$ws.Range("K31").Value = 9232.75
$ws.Range("M31").Value = -8937.75
$ws.Range("H34").Value = 433127.62
$ws.Range("I34").Value = 9232.75
$ws.Range("K34").Value = 9232.75
$ws.Range("M34").Value = -9030.75
$ws.Range("H58").Value = 2439.3684
$ws.Range("I58").Value = 2815
$ws.Range("K58").Value = 2815
$ws.Range("M58").Value = -2612
$ws.Range("H99").Value = 31376446
$ws.Range("I99").Value = 12504139
$ws.Range("K99").Value = 12504139
$ws.Range("M99").Value = -12502641
$ws.Range("H113").Value = 1878.1
$ws.Range("I113").Value = 1624.3334
$ws.Range("K113").Value = 1624.3334
$ws.Range("M113").Value = 545.6666
$ws.Range("H122").Value = 2323.5454
$ws.Range("I122").Value = 2261
$ws.Range("J122").Value = 2536.2
$ws.Range("K122").Value = 6783
$ws.Range("L122").Value = 7608.599999999999
$ws.Range("M122").Value = -4333
$ws.Range("N122").Value = -12508.6
$ws.Range("H126").Value = 31376446
$ws.Range("I126").Value = 12504139
$ws.Range("K126").Value = 37512417
$ws.Range("M126").Value = -37509947
$ws.Range("H128").Value = 89780
$ws.Range("J128").Value = 89780
$ws.Range("L128").Value = 89780
$ws.Range("N128").Value = -99740
$ws.Range("H134").Value = 1241.5883
$ws.Range("I134").Value = 1155.25
$ws.Range("K134").Value = 3465.75
$ws.Range("M134").Value = -930.75
$ws.Range("H136").Value = 2439.3684
$ws.Range("I136").Value = 2815
$ws.Range("K136").Value = 8445
$ws.Range("M136").Value = -5895
$ws.Range("H141").Value = 351844.7
$ws.Range("J141").Value = 351844.7
$ws.Range("L141").Value = 351844.7
$ws.Range("N141").Value = -362204.7

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1246.4445
$ws.Range("I5").Value = 1031.1428
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 3093.4284
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -2981.4284
$ws.Range("N5").Value = -6224
$ws.Range("H52").Value = 3912.8572
$ws.Range("J52").Value = 3912.8572
$ws.Range("L52").Value = 11738.5716
$ws.Range("N52").Value = -12270.5716
$ws.Range("H61").Value = 57
$ws.Range("I61").Value = 57
$ws.Range("K61").Value = 171
$ws.Range("M61").Value = 44
$ws.Range("H68").Value = 1704.6666
$ws.Range("I68").Value = 1244
$ws.Range("K68").Value = 3732
$ws.Range("M68").Value = -2921
$ws.Range("H71").Value = 1704.6666
$ws.Range("I71").Value = 1244
$ws.Range("K71").Value = 11196
$ws.Range("M71").Value = -7140
$ws.Range("H106").Value = 9954.333000000001
$ws.Range("I106").Value = 8996.666999999999
$ws.Range("J106").Value = 10433.167
$ws.Range("K106").Value = 26990.001
$ws.Range("L106").Value = 31299.501
$ws.Range("M106").Value = -26044.001
$ws.Range("N106").Value = -33191.501
$ws.Range("H131").Value = 4989.148
$ws.Range("I131").Value = 2423
$ws.Range("J131").Value = 5722.3335
$ws.Range("K131").Value = 7269
$ws.Range("L131").Value = 17167.0005
$ws.Range("M131").Value = -2229
$ws.Range("N131").Value = -27247.0005
$ws.Range("H135").Value = 1246.4445
$ws.Range("I135").Value = 1031.1428
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 9280.2852
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -6745.2852
$ws.Range("N135").Value = -23070
$ws.Range("H137").Value = 5306.5835
$ws.Range("I137").Value = 3097.4285
$ws.Range("K137").Value = 9292.2855
$ws.Range("M137").Value = -4192.2855

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 77000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H30").Value = 77000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H97").Value = 487.3
$ws.Range("I97").Value = 531.86664
$ws.Range("K97").Value = 531.86664
$ws.Range("M97").Value = -35.86663999999996
$ws.Range("H102").Value = 6042.2
$ws.Range("I102").Value = 6998.3335
$ws.Range("J102").Value = 5632.4287
$ws.Range("K102").Value = 6998.3335
$ws.Range("L102").Value = 5632.4287
$ws.Range("M102").Value = -5376.3335
$ws.Range("N102").Value = -8876.4287
$ws.Range("H113").Value = 2527.5293
$ws.Range("I113").Value = 2355.5557
$ws.Range("J113").Value = 2721
$ws.Range("K113").Value = 2355.5557
$ws.Range("L113").Value = 2721
$ws.Range("M113").Value = -185.5556999999999
$ws.Range("N113").Value = -7061
$ws.Range("H126").Value = 5752.048
$ws.Range("I126").Value = 5952.2
$ws.Range("K126").Value = 17856.6
$ws.Range("M126").Value = -15386.6
$ws.Range("H132").Value = 4374.7407
$ws.Range("I132").Value = 4374.7407
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13124.2221
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10594.2221
$ws.Range("N132").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 102499.5
$ws.Range("J11").Value = 135999
$ws.Range("L11").Value = 135999
$ws.Range("N11").Value = -136279
$ws.Range("H22").Value = 2667.0715
$ws.Range("I22").Value = 1219.5
$ws.Range("K22").Value = 1219.5
$ws.Range("M22").Value = -924.5
$ws.Range("H25").Value = 34200
$ws.Range("I25").Value = 35600
$ws.Range("K25").Value = 35600
$ws.Range("M25").Value = -35370
$ws.Range("H27").Value = 2667.0715
$ws.Range("I27").Value = 1219.5
$ws.Range("K27").Value = 1219.5
$ws.Range("M27").Value = -1112.5
$ws.Range("H40").Value = 33338584
$ws.Range("J40").Value = 3468.3333
$ws.Range("L40").Value = 3468.3333
$ws.Range("N40").Value = -3740.3333
$ws.Range("H46").Value = 13889996
$ws.Range("I46").Value = 16667795
$ws.Range("K46").Value = 16667795
$ws.Range("M46").Value = -16667607
$ws.Range("H55").Value = 4112.125
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -827
$ws.Range("H61").Value = 3393.4333
$ws.Range("I61").Value = 3553.7144
$ws.Range("K61").Value = 3553.7144
$ws.Range("M61").Value = -3351.7144
$ws.Range("H100").Value = 2575.7334
$ws.Range("I100").Value = 2658.9092
$ws.Range("K100").Value = 2658.9092
$ws.Range("M100").Value = -2117.9092
$ws.Range("H113").Value = 3393.4333
$ws.Range("I113").Value = 3553.7144
$ws.Range("K113").Value = 3553.7144
$ws.Range("M113").Value = -1383.7144
$ws.Range("H122").Value = 11072.182
$ws.Range("I122").Value = 6499
$ws.Range("K122").Value = 19497
$ws.Range("M122").Value = -17047
$ws.Range("H132").Value = 8638.348
$ws.Range("I132").Value = 8746.789000000001
$ws.Range("J132").Value = 8123.25
$ws.Range("K132").Value = 26240.367
$ws.Range("L132").Value = 24369.75
$ws.Range("M132").Value = -23710.367
$ws.Range("N132").Value = -29429.75
$ws.Range("H136").Value = 6938.1577
$ws.Range("I136").Value = 9143.799999999999
$ws.Range("J136").Value = 4487.4443
$ws.Range("K136").Value = 27431.4
$ws.Range("L136").Value = 13462.3329
$ws.Range("M136").Value = -24881.4
$ws.Range("N136").Value = -18562.3329
$ws.Range("H139").Value = 65000
$ws.Range("I139").Value = 45000
$ws.Range("K139").Value = 45000
$ws.Range("M139").Value = -39860

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 518.28125
$ws.Range("I113").Value = 497.92
$ws.Range("K113").Value = 1493.76
$ws.Range("M113").Value = 676.24
$ws.Range("H126").Value = 53339732
$ws.Range("I126").Value = 19615474
$ws.Range("J126").Value = 125003784
$ws.Range("K126").Value = 58846422
$ws.Range("L126").Value = 375011352
$ws.Range("M126").Value = -58843952
$ws.Range("N126").Value = -375016292
$ws.Range("H132").Value = 2323.2778
$ws.Range("I132").Value = 2262.6858
$ws.Range("J132").Value = 4444
$ws.Range("K132").Value = 6788.057400000001
$ws.Range("L132").Value = 13332
$ws.Range("M132").Value = -4258.057400000001
$ws.Range("N132").Value = -18392
$ws.Range("H135").Value = 209911.64
$ws.Range("J135").Value = 209911.64
$ws.Range("L135").Value = 209911.64
$ws.Range("N135").Value = -220051.64
